$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.273.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.651.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.509'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.256'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.86%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0628'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.882.43'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.649.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("E14").Value = '  -0.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.542'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.89'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.254.74'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.15%  '

$ws.Range("E18").Value = '  +0.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.51%  '

$ws.Range("E20").Value = '  -0.53%  '

$ws.Range("E21").Value = '  +2.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.08%  '

$ws.Range("E24").Value = '  +0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.13%  '

$ws.Range("E26").Value = '  -0.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.22%  '

$ws.Range("E28").Value = '  +1.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.89'
$ws.Range("D29").Style = "Normal"

$ws.Range("E30").Value = '  -0.62%  '

$ws.Range("E32").Value = '  -0.32%  '

$ws.Range("E33").Value = '  +1.08%  '

$ws.Range("E34").Value = '  +1.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.274.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.05%  '

$ws.Range("E37").Value = '  +3.23%  '

$ws.Range("E38").Value = '  +2.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.845'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.810'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("E42").Value = '  +0.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.19'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.791.92'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("E48").Value = '  +9.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0514'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.73'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0978'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.76%  '
